$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new log row (row 44) mirroring the existing "Object Detection"
# entries for end device 2, with a new timestamp.
$row = 44
$ws.Cells.Item($row, 1).Value = $row
$ws.Cells.Item($row, 2).Value = "Object Detection"
$ws.Cells.Item($row, 3).Value = "New image from end device 2 detected and recorded to database."
$ws.Cells.Item($row, 4).Value = "07/05/2022 03:36:31 AM"
